# Rerun and summarise models without urban landuse:
#  - rename the 9 summary sheets to their new summ<id> names
#  - update the "Education[T.Unknown]" label to "Education[T.Unknown/Other]"
#    on every sheet (row 5, column A)

$wb = $excel.ActiveWorkbook

$oldNames = @(
    "summ27253330",
    "summ27581494",
    "summ27876547",
    "summ28144693",
    "summ28426993",
    "summ28782025",
    "summ29076002",
    "summ29345259",
    "summ29587055"
)

$newNames = @(
    "summ07664943",
    "summ07879573",
    "summ08116258",
    "summ08376375",
    "summ08622422",
    "summ08883081",
    "summ09152479",
    "summ09426045",
    "summ09671568"
)

for ($i = 0; $i -lt $oldNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($oldNames[$i])
    $ws.Name = $newNames[$i]
}

foreach ($ws in $wb.Worksheets) {
    if ($ws.Range("A5").Value2 -eq "Education[T.Unknown]") {
        $ws.Range("A5").Value = "Education[T.Unknown/Other]"
    }
}
